$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.510.78"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.645.99"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.65"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.530"
$ws.Range("E6").Value = "  +3.92%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.56"
$ws.Range("E8").Value = "  -2.22%  "

$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.48"
$ws.Range("E12").Value = "  -1.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.645.23"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.589"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.461.26"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.51"
$ws.Range("E18").Value = "  -4.28%  "

$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -1.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  -3.64%  "

$ws.Range("E23").Value = "  +3.76%  "

$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.81"
$ws.Range("E25").Value = "  +0.32%  "

$ws.Range("E26").Value = "  -2.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.114"
$ws.Range("E27").Value = "  +1.65%  "

$ws.Range("E28").Value = "  -4.45%  "

$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  -3.48%  "

$ws.Range("E31").Value = "  -3.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.426.30"
$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("E35").Value = "  +0.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.569"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.890"
$ws.Range("E38").Value = "  -4.31%  "

$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("E40").Value = "  -1.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.823"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +2.66%  "

$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.93"
$ws.Range("E46").Value = "  -7.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.788.41"
$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("E48").Value = "  -3.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.57"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0107"
$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("E51").Value = "  -3.15%  "
